# "add transaction checks, single day fix, unit testing update"
#
# The "type" column (B) used the mixed-case label "Buy" for every
# transaction row; normalize it to the upper-case "BUY" used by the
# new transaction-checking logic. Also restore the active selection
# to Q9 (previously the whole column G was selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row   # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 2).Value2 -eq "Buy") {
        $ws.Cells.Item($r, 2).Value = "BUY"
    }
}

$ws.Range("Q9").Select()
